$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = "Qosimova Nigoraxon Mahmudjon qizi"; B = "Maktabgacha talim tashkiloti metodisti"; C = "AB7943865"; D = "262"; E = $null; F = "Paxtaobod tumani"; G = "998990335588"; H = "20-11-2024" },
    @{ A = "Dadaboyeva Marhabo G'aniyevna"; B = "Maktabgacha talim tashkiloti psixologi"; C = "AB4871365"; D = "263"; E = $null; F = "Paxtaobod tumani"; G = "998994314929"; H = "20-11-2024" },
    @{ A = "Jumaboyeva Mohiyaxon Ismoilovna"; B = "Maktabgacha talim tashkiloti direktori"; C = "AD2021085"; D = "264"; E = $null; F = "Paxtaobod tumani"; G = "998999055045"; H = "20-11-2024" },
    @{ A = "Yunusova Gulnoza Xikmatilla qizi"; B = "Maktabgacha talim tashkiloti tarbiyachisi"; C = "AB9919991"; D = "265"; E = $null; F = "Yuqori Chirchiq tumani"; G = "998943789199"; H = "20-11-2024" },
    @{ A = "Soxibova Shaxodat Komilovna"; B = "Maktabgacha talim tashkiloti direktori"; C = "AB0831483"; D = "266"; E = "Namangan viloyati"; F = "Uchqoʻrgʻon tumani"; G = "998974675009"; H = "20-11-2024" },
    @{ A = "Maxammadiyeva Surayyo Normuradovna"; B = "Maktabgacha talim tashkiloti metodisti"; C = "AD7998828"; D = "267"; E = "Navoiy viloyati"; F = "Navoiy shahri"; G = "998943799003"; H = "20-11-2024" }
)

$startRow = 84
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    # Numeric-looking IDs must stay text (leading apostrophe forces text, like the rest of the column).
    $ws.Cells.Item($r, 4).Value = "'" + $data.D
    if ($null -ne $data.E) {
        $ws.Cells.Item($r, 5).Value = $data.E
    }
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = "'" + $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
}
